# 🔄 Actualización automática del tracker
# Updates row 167 (result now known) and appends new pending rows 170-181.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $text) {
    # Force a literal text value (avoids Excel's automatic date/number
    # parsing for strings like "2025-09-17"), then strip the formatting
    # that gets attached by the temporary "@" number format so the cell
    # ends up with no explicit style - matching freshly authored rows.
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Row 167: result came in ---
Set-TextCell $ws "G167" "Fallo"
$ws.Range("H167").Value = -1

# --- New rows 170-181 ---
$rows = @(
    @{ Row=170; A=14688912; B="2025-09-17"; C="Zizou Bergs";              D="Juncheng Shang";        E="Gana Zizou Bergs";                 F=1.91 },
    @{ Row=171; A=14688810; B="2025-09-17"; C="Adrian Mannarino";         D="Yibing Wu";              E="Gana Adrian Mannarino";            F=1.67 },
    @{ Row=172; A=14655423; B="2025-09-17"; C="Clara Tauson";             D="Eva Lys";                E="Gana Clara Tauson";                F=1.53 },
    @{ Row=173; A=14655426; B="2025-09-17"; C="Katerina Siniakova";       D="Daria Kasatkina";        E="Gana Daria Kasatkina";             F=1.73 },
    @{ Row=174; A=14690566; B="2025-09-16"; C="Trey Hilderbrand";         D="Hugo Grenier";           E="Gana Trey Hilderbrand";            F=3.25 },
    @{ Row=175; A=14693670; B="2025-09-16"; C="Jack Anthrop";             D="Antoine Ghibaudo";       E="Gana Jack Anthrop";                F=2.75 },
    @{ Row=176; A=14681379; B="2025-09-16"; C="Tyler Zink";               D="Alexander Bernard";      E="Gana Alexander Bernard";           F=2.75 },
    @{ Row=177; A=14693672; B="2025-09-16"; C="Andre Ilagan";             D="Sebastian Dominko";      E="Gana Sebastian Dominko";           F=3.25 },
    @{ Row=178; A=14681382; B="2025-09-16"; C="Martin Damm Jr";           D="Aidan Kim";              E="Gana Aidan Kim";                   F=3 },
    @{ Row=179; A=14692728; B="2025-09-16"; C="Juan Pablo Varillas";      D="Santiago De la Fuente";  E="Gana Santiago De la Fuente";       F=3.75 },
    @{ Row=180; A=14692731; B="2025-09-16"; C="Luciano Emanuel Ambrogi";  D="Facundo Bagnis";         E="Gana Luciano Emanuel Ambrogi";     F=2.75 },
    @{ Row=181; A=14692727; B="2025-09-16"; C="Matias Soto";              D="Carlos Maria Zarate";    E="Gana Carlos Maria Zarate";         F=4.5 }
)

foreach ($entry in $rows) {
    $rowNum = $entry.Row
    $ws.Range("A$rowNum").Value = $entry.A
    Set-TextCell $ws "B$rowNum" $entry.B
    Set-TextCell $ws "C$rowNum" $entry.C
    Set-TextCell $ws "D$rowNum" $entry.D
    Set-TextCell $ws "E$rowNum" $entry.E
    $ws.Range("F$rowNum").Value = $entry.F
}
